$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '19.890.08'

# Row 3
$ws.Range("D3").Value = '1.392.91'
$ws.Range("E3").Value = '  -8.98%  '

# Row 4
$style_D4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("D4").Style = $style_D4
$ws.Range("E4").Value = '  -0.17%  '

# Row 5
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9990'
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  -0.22%  '

# Row 6
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.51'
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  -6.72%  '

# Row 7
$style_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3648'
$ws.Range("D7").Style = $style_D7
$ws.Range("E7").Value = '  -7.77%  '

# Row 8
$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3070'
$ws.Range("D8").Style = $style_D8
$ws.Range("E8").Value = '  -2.71%  '

# Row 9
$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.23'
$ws.Range("D9").Style = $style_D9
$ws.Range("E9").Value = '  -6.71%  '

# Row 10
$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9906'
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  -5.69%  '

# Row 11
$style_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06423'
$ws.Range("D11").Style = $style_D11
$ws.Range("E11").Value = '  -10.09%  '

# Row 12
$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9992'
$ws.Range("D12").Style = $style_D12
$ws.Range("E12").Value = '  -0.18%  '

# Row 13
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.339'
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = '  -6.01%  '

# Row 14
$ws.Range("E14").Value = '  -7.06%  '

# Row 15
$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.083'
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = '  -7.69%  '

# Row 16
$ws.Range("D16").Value = '1.390.68'
$ws.Range("E16").Value = '  -9.68%  '

# Row 17
$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009974'
$ws.Range("D17").Style = $style_D17

# Row 18
$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05645'
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  -14.41%  '

# Row 19
$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9989'
$ws.Range("D19").Style = $style_D19
$ws.Range("E19").Value = '  -0.15%  '

# Row 20
$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.27'
$ws.Range("D20").Style = $style_D20
$ws.Range("E20").Value = '  -16.01%  '

# Row 21
$ws.Range("E21").Value = '  -9.79%  '

# Row 22
$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.64'
$ws.Range("D22").Style = $style_D22
$ws.Range("E22").Value = '  -5.22%  '

# Row 23
$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.80'
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  +1.61%  '

# Row 24
$style_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.243'
$ws.Range("D24").Style = $style_D24
$ws.Range("E24").Value = '  -4.95%  '

# Row 25
$ws.Range("D25").Value = '19.893.78'
$ws.Range("E25").Value = '  -7.95%  '

# Row 26
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.180'
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = '  -6.17%  '

# Row 27
$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '135.20'
$ws.Range("D27").Style = $style_D27
$ws.Range("E27").Value = '  -9.68%  '

# Row 28
$style_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.62'
$ws.Range("D28").Style = $style_D28
$ws.Range("E28").Value = '  -9.05%  '

# Row 29
$ws.Range("D29").Value = '1.548.48'
$ws.Range("E29").Value = '  -9.44%  '

# Row 30
$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '108.38'
$ws.Range("D30").Style = $style_D30
$ws.Range("E30").Value = '  -7.07%  '

# Row 31
$style_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.027'
$ws.Range("D31").Style = $style_D31
$ws.Range("E31").Value = '  -16.91%  '

# Row 32
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.212'
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  -13.91%  '

# Row 33
$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7974'
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = '  -14.36%  '

# Row 34
$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07588'
$ws.Range("D34").Style = $style_D34
$ws.Range("E34").Value = '  -6.51%  '

# Row 35
$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.276'
$ws.Range("D35").Style = $style_D35
$ws.Range("E35").Value = '  -1.16%  '

# Row 36
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9987'
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  -0.17%  '

# Row 37
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05653'
$ws.Range("D37").Style = $style_D37
$ws.Range("E37").Value = '  -5.18%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.714'
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = '  -8.10%  '

# Row 39
$ws.Range("B39").Value = 'WEMIXTOKEN'
$ws.Range("C39").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.367'
$ws.Range("D39").Style = $style_D39
$ws.Range("E39").Value = '  -5.92%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02036'
$ws.Range("D40").Style = $style_D40
$ws.Range("E40").Value = '  -7.49%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1893'
$ws.Range("D41").Style = $style_D41
$ws.Range("E41").Value = '  -6.47%  '

# Row 42
$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.12'
$ws.Range("D42").Style = $style_D42
$ws.Range("E42").Value = '  -7.74%  '

# Row 43
$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.063'
$ws.Range("D43").Style = $style_D43
$ws.Range("E43").Value = '  -9.65%  '

# Row 44
$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5203'
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = '  -9.65%  '

# Row 45
$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.474'
$ws.Range("D45").Style = $style_D45
$ws.Range("E45").Value = '  -6.41%  '

# Row 46
$style_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.02'
$ws.Range("D46").Style = $style_D46
$ws.Range("E46").Value = '  -6.98%  '

# Row 47
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4997'
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  -8.85%  '

# Row 48
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '109.91'
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = '  -4.80%  '

# Row 49
$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.742'
$ws.Range("D49").Style = $style_D49
$ws.Range("E49").Value = '  -6.66%  '

# Row 50
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9990'
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = '  -0.15%  '

# Row 51
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.025'
$ws.Range("D51").Style = $style_D51
$ws.Range("E51").Value = '  -11.92%  '
